# Updated cryptos list on Sat Mar 25 23:33:51 UTC 2023 with GitHub Actions
#
# Refresh the "Price" (D) and "Volume(1h)" (E) columns for each coin row
# with the latest scraped values. Two rows (39/40) also swap which coin
# (TheSandbox / InternetComputer(DFINITY)) occupies that ranking slot, so
# their Coin name (B) and Link (C) columns are updated too.
#
# Price values that look like plain numbers ("322.60", "1.002", ...) are
# written with a leading single-quote so Excel keeps them as literal text
# (preserving trailing zeros etc.) instead of silently converting them to
# numbers, exactly like the original cells (stored as inline text strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.502.69'
$ws.Range("E2").Value = '  +0.20%  '

$ws.Range("D3").Value = '1.742.82'
$ws.Range("E3").Value = '  -0.43%  '

$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").Value = '''322.60'
$ws.Range("E5").Value = '  +0.20%  '

$ws.Range("D6").Value = '''1.002'
$ws.Range("E6").Value = '  +0.27%  '

$ws.Range("D7").Value = '''0.4432'
$ws.Range("E7").Value = '  +4.41%  '

$ws.Range("D8").Value = '''0.3520'
$ws.Range("E8").Value = '  -2.30%  '

$ws.Range("D9").Value = '''0.07409'
$ws.Range("E9").Value = '  -0.79%  '

$ws.Range("D10").Value = '''41.56'
$ws.Range("E10").Value = '  -1.87%  '

$ws.Range("D11").Value = '''1.075'
$ws.Range("E11").Value = '  -2.18%  '

$ws.Range("D12").Value = '''1.002'
$ws.Range("E12").Value = '  +0.18%  '

$ws.Range("D13").Value = '''20.42'
$ws.Range("E13").Value = '  -1.25%  '

$ws.Range("D14").Value = '''5.904'
$ws.Range("E14").Value = '  -2.08%  '

$ws.Range("D15").Value = '''7.072'
$ws.Range("E15").Value = '  -2.16%  '

$ws.Range("D16").Value = '1.741.53'
$ws.Range("E16").Value = '  -0.13%  '

$ws.Range("D17").Value = '''91.38'
$ws.Range("E17").Value = '  -1.62%  '

$ws.Range("D18").Value = '''0.00001051'
$ws.Range("E18").Value = '  -1.12%  '

$ws.Range("D19").Value = '''0.06379'
$ws.Range("E19").Value = '  -0.02%  '

$ws.Range("E20").Value = '  +0.25%  '

$ws.Range("D21").Value = '''16.81'
$ws.Range("E21").Value = '  -1.33%  '

$ws.Range("D22").Value = '''5.714'
$ws.Range("E22").Value = '  -2.92%  '

$ws.Range("D23").Value = '27.539.26'
$ws.Range("E23").Value = '  +0.10%  '

$ws.Range("D24").Value = '''11.11'
$ws.Range("E24").Value = '  -1.06%  '

$ws.Range("E25").Value = '  -0.16%  '

$ws.Range("D26").Value = '''160.33'
$ws.Range("E26").Value = '  -0.39%  '

$ws.Range("D27").Value = '''20.02'
$ws.Range("E27").Value = '  -1.22%  '

$ws.Range("D28").Value = '1.941.46'
$ws.Range("E28").Value = '  -0.09%  '

$ws.Range("E29").Value = '  +0.68%  '

$ws.Range("D30").Value = '''2.033'
$ws.Range("E30").Value = '  -4.53%  '

$ws.Range("E31").Value = '  -5.33%  '

$ws.Range("D32").Value = '''0.09082'
$ws.Range("E32").Value = '  +2.53%  '

$ws.Range("E33").Value = '  -0.07%  '

$ws.Range("D34").Value = '''5.359'
$ws.Range("E34").Value = '  -3.34%  '

$ws.Range("D35").Value = '''0.02269'
$ws.Range("E35").Value = '  -0.69%  '

$ws.Range("D36").Value = '''11.59'
$ws.Range("E36").Value = '  -5.12%  '

$ws.Range("D37").Value = '''0.06026'
$ws.Range("E37").Value = '  +0.52%  '

$ws.Range("D38").Value = '''0.2060'
$ws.Range("E38").Value = '  -1.61%  '

$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = '''4.886'
$ws.Range("E39").Value = '  -1.00%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '''0.6225'
$ws.Range("E40").Value = '  -1.40%  '

$ws.Range("D41").Value = '''1.185'
$ws.Range("E41").Value = '  +0.27%  '

$ws.Range("E42").Value = '  -1.03%  '

$ws.Range("D43").Value = '''7.712'
$ws.Range("E43").Value = '  -2.21%  '

$ws.Range("D44").Value = '''13.19'
$ws.Range("E44").Value = '  -1.69%  '

$ws.Range("D45").Value = '''3.700'
$ws.Range("E45").Value = '  +0.33%  '

$ws.Range("D46").Value = '''0.5783'
$ws.Range("E46").Value = '  -1.58%  '

$ws.Range("D47").Value = '''121.83'

$ws.Range("D48").Value = '''1.922'
$ws.Range("E48").Value = '  -2.15%  '

$ws.Range("D49").Value = '''0.06837'
$ws.Range("E49").Value = '  +0.17%  '

$ws.Range("E50").Value = '  -4.68%  '

$ws.Range("D51").Value = '''71.32'
$ws.Range("E51").Value = '  -2.36%  '
